# LH_TC_NOTIFICATION_REVIEWS - v2.1 reviewed and closed reviewer verification
# for notification testcases (LH_TC_NOTIFICATION_REVIEW)

$wb = $excel.ActiveWorkbook

# --- Sheet: LH_TC_NOTIFICATION_REVIEWS ---
# Fix the Review ID naming convention/typo: "LH-TC-NOTIGICATION-REVIEW-00X"
# -> "LH-REVIEW-TC-NOTIFICATION-00X"
$wsReviews = $wb.Worksheets.Item("LH_TC_NOTIFICATION_REVIEWS")

$wsReviews.Range("B5").Value = "LH-REVIEW-TC-NOTIFICATION-004"
$wsReviews.Range("B2").Value = "LH-REVIEW-TC-NOTIFICATION-001"
$wsReviews.Range("B3").Value = "LH-REVIEW-TC-NOTIFICATION-002"
$wsReviews.Range("B4").Value = "LH-REVIEW-TC-NOTIFICATION-003"

# Reflect the reviewer's selection after editing the Review ID column
# (range B2:B5, anchored/active at the last-edited cell B5)
$wsReviews.Range("B2:B5").Select()

# --- Sheet: Version History ---
# Re-affirm the authors for the latest two entries (v2.0 / V2.1) so the
# row heights recompute for the wrapped "Mahmoud Abdelmageed" text.
$wsHistory = $wb.Worksheets.Item("Version History")
$wsHistory.Range("B5").Value = "Mahmoud Abdelmageed"
$wsHistory.Range("B6").Value = "Mahmoud Abdelmageed"
